# "Commenting and removed some legacy functions"
#
# The "Rectangle 10" shape (shape id 11, creationId
# {63B4AB3E-2992-4B5F-9906-4B2720BF2F49}) on slide 1 lists DeckOfCards'
# Functions. Capitalise the "draw" bullet to "Draw" and append a new
# bulleted line "<<" below it (same bullet/indent style as its siblings).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# All of the diagram's shapes live inside one top-level group.
$grp = $s.Shapes.Item(1)

$target = $null
for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $candidate = $grp.GroupItems.Item($i)
    if ($candidate.Id -eq 11) {
        $target = $candidate
    }
}

$tr = $target.TextFrame.TextRange
$count = $tr.Paragraphs().Count

# Re-case the existing run in place (touching only the single run so the
# XML keeps one <a:r> instead of being split into several runs).
$drawParagraph = $tr.Paragraphs($count, 1)
$drawParagraph.Runs(1).Text = "Draw"

# Append a new paragraph after "Draw" carrying the same bullet formatting
# (it inherits pPr/bullet/run formatting from the paragraph it follows).
$tr.InsertAfter("`r<<") | Out-Null
